$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S (year 2022) mirrors column R's formatting.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").NumberFormat = "0.0"
$ws.Range("S5").Value = 42

# Match the saved selection from the source workbook.
$ws.Range("U4").Select()
